# Applies the AfDD 2023 Annex Table Tab06 update:
#  - Refreshes the "most recent measure" year range in the header row
#    from 2012-21 to 2013-22 for the 8 literacy-rate indicator columns.
#  - Updates the PALOP / MERCOSUR footnote text (as captured verbatim
#    in the authoritative diff, including its mangled accented letters).
#  - Refreshes a batch of recalculated aggregate figures (World, Latin
#    America & Caribbean, World outside Africa, ROW income groupings,
#    Africa LDCs, Small Island Developing States) in columns C:J.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab06")

# ---------------------------------------------------------------
# 1. Header row (row 2): bump the reference period 2012-21 -> 2013-22
# ---------------------------------------------------------------
$ws.Range("C2").Value = "Youth literacy rate, population 15-24 years, both sexes (%, most recent measure 2013-22)"
$ws.Range("D2").Value = "Youth literacy rate, population 15-24 years, female (%, most recent measure 2013-22)"
$ws.Range("E2").Value = "Youth literacy rate, population 15-24 years, male (%, most recent measure 2013-22)"
$ws.Range("F2").Value = "Youth literacy rate, population 15-24 years, adjusted gender parity index (GPIA, most recent measure 2013-22)"
$ws.Range("G2").Value = "Adult literacy rate, population 15+ years, both sexes (%, most recent measure 2013-22)"
$ws.Range("H2").Value = "Adult literacy rate, population 15+ years, female (%, most recent measure 2013-22)"
$ws.Range("I2").Value = "Adult literacy rate, population 15+ years, male (%, most recent measure 2013-22)"
$ws.Range("J2").Value = "Adult literacy rate, population 15+ years, adjusted gender parity index (GPIA, most recent measure 2013-22)"

# ---------------------------------------------------------------
# 2. Footnote text update (cell A103)
# ---------------------------------------------------------------
$ws.Range("A103").Value = 'Regional Economic Communities:CEN-SAD = "Community of Sahel-Saharan States";COMESA = "Common Market for Eastern and Southern Africa";EAC = "East African Community";ECCAS = "Economic Community of Central African States";ECOWAS = "Economic Community of West African States";IGAD = "Intergovernmental Authority on Development";SADC = "Southern African Development Community";UMA = "Arab Maghreb Union";PALOP = "Pa>ses Africanos de L>ngua Oficial Portuguesa";ASEAN = "Association of Southeast Asian Nations";MERCOSUR = "Mercado Com>n del Sur".EU27 = "European Union (27 members)".OECD = "Organisation for Economic Co-operation and Development".'

# ---------------------------------------------------------------
# 3. Recalculated aggregate values (columns C:J)
# ---------------------------------------------------------------

# Row 63 - World outside Africa
$ws.Range("C63").Value = 97.124107380952395
$ws.Range("D63").Value = 96.993711666666698
$ws.Range("E63").Value = 97.284585357142902
$ws.Range("F63").Value = 0.99524880952381001
$ws.Range("G63").Value = 92.493392
$ws.Range("H63").Value = 91.049872470588198
$ws.Range("I63").Value = 93.952597294117595
$ws.Range("J63").Value = 0.96285776470588003

# Row 64 - Latin America and Caribbean
$ws.Range("C64").Value = 97.5552858333334
$ws.Range("D64").Value = 97.856585416666704
$ws.Range("E64").Value = 97.264937083333393
$ws.Range("F64").Value = 1.0060150000000001
$ws.Range("G64").Value = 92.2380876
$ws.Range("H64").Value = 91.667259200000004
$ws.Range("I64").Value = 92.8557536
$ws.Range("J64").Value = 0.98584000000000005

# Row 65 - Asia (no high inc.) -- only D65 has a (negligible, last-digit) update
$ws.Range("D65").Value = 93.293322222222301

# Row 66 - World
$ws.Range("C66").Value = 90.209619402985098
$ws.Range("D66").Value = 89.180483507462696
$ws.Range("E66").Value = 91.359260746268703
$ws.Range("F66").Value = 0.96591686567164003
$ws.Range("G66").Value = 83.185730592592606
$ws.Range("H66").Value = 79.986195851851903
$ws.Range("I66").Value = 86.552861185185193
$ws.Range("J66").Value = 0.90171096296296005

# Row 83 - ROW, Non-resource-rich countries
$ws.Range("C83").Value = 96.951102121212202
$ws.Range("D83").Value = 96.736173333333397
$ws.Range("E83").Value = 97.182109393939498
$ws.Range("F83").Value = 0.99314181818182001
$ws.Range("G83").Value = 91.962249999999997
$ws.Range("H83").Value = 90.476356119402993
$ws.Range("I83").Value = 93.4991217910448
$ws.Range("J83").Value = 0.96030671641790999

# Row 87 - ROW, Lower middle income countries
$ws.Range("C87").Value = 94.573083999999994
$ws.Range("D87").Value = 94.368003999999999
$ws.Range("E87").Value = 94.762525600000004
$ws.Range("F87").Value = 0.99463760000000001
$ws.Range("G87").Value = 85.179221200000001
$ws.Range("H87").Value = 82.177672000000001
$ws.Range("I87").Value = 88.281700400000005
$ws.Range("J87").Value = 0.92352319999999999

# Row 89 - ROW, Upper middle income countries
$ws.Range("C89").Value = 98.864615000000001
$ws.Range("D89").Value = 98.959063157894803
$ws.Range("E89").Value = 98.776679736842098
$ws.Range("F89").Value = 1.0017989473684199
$ws.Range("G89").Value = 95.982192368420996
$ws.Range("H89").Value = 95.313946315789494
$ws.Range("I89").Value = 96.665577368421097
$ws.Range("J89").Value = 0.98567657894737004

# Row 91 - Africa, Least Developed Countries -- only G91 has a negligible update
$ws.Range("G91").Value = 57.541360322580701

# Row 94 - ROW, Small Island Developing States
$ws.Range("C94").Value = 96.3071985714286
$ws.Range("D94").Value = 96.683132857142894
$ws.Range("E94").Value = 95.948515
$ws.Range("F94").Value = 1.00761428571429
$ws.Range("G94").Value = 91.108957333333393
$ws.Range("H94").Value = 90.690535999999994
$ws.Range("I94").Value = 91.531800666666697
$ws.Range("J94").Value = 0.98721800000000004
